$d = $word.ActiveDocument

function FindParaIndexByStart($doc, $s) {
    $count = $doc.Paragraphs.Count
    for ($i = 1; $i -le $count; $i++) {
        $p = $doc.Paragraphs.Item($i)
        if ($p.Range.Start -eq $s) {
            return $i
        }
    }
    return -1
}

# ---------------------------------------------------------------------
# 1) Merge the split "...nevěří..." runs (removes the stray proofErr
#    gramStart/gramEnd wrapper) back into a single run of text.
# ---------------------------------------------------------------------
$find1 = $d.Content.Find
$find1.Execute(", který zpočátku nevěří, že je", $true, $false, $false, $false, $false, $true, 1, $false, ", který zpočátku nevěří, že je", 2) | Out-Null

# ---------------------------------------------------------------------
# 2) Merge the split "...snaží..." runs the same way.
# ---------------------------------------------------------------------
$find2 = $d.Content.Find
$find2.Execute("se velmi snaží, je šikovná", $true, $false, $false, $false, $false, $true, 1, $false, "se velmi snaží, je šikovná", 2) | Out-Null

# ---------------------------------------------------------------------
# 3) Insert a new "Romain Rolland" sub-heading (with two bullet points)
#    just above the "Německý" heading, and retitle that heading
#    "Německo".
# ---------------------------------------------------------------------
$findHeading = $d.Content
$findHeading.Find.Execute("Německý") | Out-Null
$headingStart = $findHeading.Start
$headingIdx = FindParaIndexByStart $d $headingStart

# Make room: three fresh paragraphs ahead of the "Německý" heading.
$insertPoint = $d.Range($headingStart, $headingStart)
$insertPoint.InsertParagraphBefore() | Out-Null
$insertPoint.InsertParagraphBefore() | Out-Null
$insertPoint.InsertParagraphBefore() | Out-Null

# --- Paragraph 1: "Romain Rolland" heading (style Nadpis4) ---
$pRolland = $d.Paragraphs.Item($headingIdx)
$pRolland.Style = "Nadpis 4"
$pRolland.Range.Text = "Romain Rolland"
$rollandStart = $pRolland.Range.Start
# "Romain " (non-italic, matches surrounding heading 4 author-name run)
$d.Range($rollandStart, $rollandStart + 7).Font.Italic = 0
# "Rolland" (non-italic, bold, underlined - the surname run)
$surnameRange = $d.Range($rollandStart + 7, $rollandStart + 15)
$surnameRange.Font.Italic = 0
$surnameRange.Font.Bold = 1
$surnameRange.Font.Underline = 1

# --- Paragraph 2: bullet "Prozaik, dramatik" ---
$pBullet1 = $d.Paragraphs.Item($headingIdx + 1)
$pBullet1.Style = "Odstavecseseznamem"
$pBullet1.Range.Text = "Prozaik, dramatik"
$pBullet1.Range.ListFormat.ApplyBulletDefault() | Out-Null

# --- Paragraph 3: bullet "Petr a Lucie" ---
$pBullet2 = $d.Paragraphs.Item($headingIdx + 2)
$pBullet2.Style = "Odstavecseseznamem"
$pBullet2.Range.Text = "Petr a Lucie"
$pBullet2.Range.ListFormat.ApplyBulletDefault() | Out-Null

# --- Retitle the (now shifted) "Německý" heading to "Německo" ---
$findNemecky = $d.Content
$findNemecky.Find.Execute("Německý", $true, $false, $false, $false, $false, $true, 1, $false, "Německo", 2) | Out-Null

# ---------------------------------------------------------------------
# 4) Retitle the "Anglický" heading to "Anglie".
# ---------------------------------------------------------------------
$findAnglicky = $d.Content
$findAnglicky.Find.Execute("Anglický", $true, $false, $false, $false, $false, $true, 1, $false, "Anglie", 2) | Out-Null

Write-Output "done"
